$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $range = $d.Content
    $range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-Text "2025-11-02 Sunday" "2025-11-03 Monday"

Replace-Text "87×89=" "12×60="
Replace-Text "66×55=" "98×14="
Replace-Text "54×24=" "13×40="
Replace-Text "72×74=" "37×84="
Replace-Text "29×57=" "16×95="

Replace-Text "74×44=" "35×41="
Replace-Text "36×49=" "71×99="
Replace-Text "41×49=" "90×46="
Replace-Text "70×11=" "83×38="
Replace-Text "44×41=" "44×61="

Replace-Text "33×97=" "34×16="
Replace-Text "93×37=" "80×83="
Replace-Text "44×65=" "59×56="
Replace-Text "59×22=" "66×73="
Replace-Text "23×84=" "64×60="

Replace-Text "84×98=" "41×80="
Replace-Text "50×85=" "95×57="
Replace-Text "66×85=" "90×30="
Replace-Text "74×81=" "62×56="
Replace-Text "73×44=" "18×27="

Replace-Text "76×88=" "59×28="
Replace-Text "89×85=" "13×46="
Replace-Text "46×81=" "65×11="
Replace-Text "52×16=" "90×85="
Replace-Text "96×45=" "57×96="
